$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.507.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.741.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4454"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3525"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.075"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.892"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.065"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.738.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06370"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.719"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.542.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.097"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.943.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "124.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.026"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.042"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09042"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.654"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.356"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02267"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.86%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2057"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6229"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.870"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.183"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.376"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.705"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5781"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.921"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06837"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.109"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.63%  "
